$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title shape: "Testing" + " " + "custom" + " " + "properties"
# -> "Testing " + "custom " + "properties"
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Characters(1, 8).Text = "Testing "
$titleRange.Characters(9, 7).Text = "custom "

# Subtitle shape: two line breaks then "A." + " " + "M."
# -> two line breaks then "A. " + "M."
$subtitleRange = $s.Shapes.Item(2).TextFrame.TextRange
$subtitleRange.Characters(3, 3).Text = "A. "
